$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# ---------------------------------------------------------------------------
# Update the two period-header rows (row 8 = financial period labels,
# row 9 = publish dates). Every year column shifts one period forward:
# the oldest period (1396/12) is dropped and a new period (1401/12) is
# appended in column H.
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

$ws.Range("D9").Value = "1399-03-21 (8)"
$ws.Range("E9").Value = "1400-03-02 (8)"
$ws.Range("F9").Value = "1401-03-08 (8)"
$ws.Range("G9").Value = "1402-02-28 (7)"
$ws.Range("H9").Value = "1402-02-28"

# ---------------------------------------------------------------------------
# Update the balance-sheet figures. Each data row shifts left by one year
# column (D<-E, E<-F, F<-G, G<-H) and receives the new figure for the
# latest period (1401/12) in column H.
# ---------------------------------------------------------------------------
function Set-Row($r, $d, $e, $f, $g, $h) {
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
}

Set-Row 12 158361  147377  58912   282794  426092
Set-Row 14 1716486 2262059 4208929 5048473 9246434
Set-Row 15 552415  918350  2039000 2667633 4583450
Set-Row 16 232171  199390  234218  524170  1456570
Set-Row 18 2659433 3527176 6541059 8523070 15712546
Set-Row 19 2619    1944    8202    15640   32908
Set-Row 20 291269  328799  328807  374932  373195
Set-Row 22 404700  559160  669370  1223388 2256464
Set-Row 23 17607   8122    5145    2837    8055
Set-Row 26 716195  898025  1011524 1616797 2670622
Set-Row 27 3375628 4425201 7552583 10139867 18383168
Set-Row 29 339064  325275  338892  884747  2697941
Set-Row 31 0       0       0       0       7479
Set-Row 32 232246  462917  620213  607895  1035247
Set-Row 33 20559   79553   424424  363263  131361
Set-Row 34 853079  776280  1937869 3298895 6391971
Set-Row 37 1444948 1644025 3321398 5154800 10263999
Set-Row 38 0       261317  0       0       0
Set-Row 41 70903   98692   153951  238536  460259
Set-Row 42 70903   360009  153951  238536  460259
Set-Row 43 1515851 2004034 3475349 5393336 10724258
Set-Row 45 1030000 1030000 1330000 2430000 3630000
Set-Row 48 0       0       0       -1060   -7969
Set-Row 49 0       0       0       314     0
Set-Row 50 96225   103000  133000  223222  363000
Set-Row 56 733552  1288167 2614234 2094055 3673879
Set-Row 57 1859777 2421167 4077234 4746531 7658910
Set-Row 58 3375628 4425201 7552583 10139867 18383168

# Rows whose values are the literal "-" placeholder string in every year
# column stay "-" in every column (D39/D52/D54 flip from 0 to "-"; D49
# flips from "-" to 0).
$ws.Range("D39").Value = "-"
$ws.Range("D49").Value = 0
$ws.Range("D52").Value = "-"
$ws.Range("D54").Value = "-"
